# Refresh the "Price" / "Volume(1h)" columns of the crypto table with
# the latest GitHub-Actions scrape, and swap the #50 listing from
# Stacks to Aave.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.203.57"
$ws.Range("E2").Value = "  +3.06%  "
$ws.Range("D3").Value = "2.248.93"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'295.04"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Value = "'87.37"
$ws.Range("E6").Value = "  +8.90%  "
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.476"
$ws.Range("E9").Value = "  +3.40%  "
$ws.Range("D10").Value = "'31.29"
$ws.Range("E10").Value = "  +12.06%  "
$ws.Range("D11").Value = "'0.0804"
$ws.Range("E11").Value = "  +4.09%  "
$ws.Range("D12").Value = "'47.25"
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("E14").Value = "  +6.03%  "
$ws.Range("D15").Value = "2.593.82"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").Value = "'14.32"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").Value = "2.241.23"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "'0.737"
$ws.Range("E18").Value = "  +2.81%  "
$ws.Range("D19").Value = "40.108.90"
$ws.Range("E19").Value = "  +3.03%  "
$ws.Range("E20").Value = "  +4.05%  "
$ws.Range("D21").Value = "'5.87"
$ws.Range("E21").Value = "  +2.51%  "
$ws.Range("D22").Value = "'10.72"
$ws.Range("E22").Value = "  +8.90%  "
$ws.Range("D23").Value = "'65.95"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").Value = "'237.30"
$ws.Range("E24").Value = "  +5.55%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  +3.67%  "
$ws.Range("E27").Value = "  +6.80%  "
$ws.Range("D28").Value = "'23.27"
$ws.Range("E28").Value = "  +4.63%  "
$ws.Range("E29").Value = "  +1.74%  "
$ws.Range("D30").Value = "'9.32"
$ws.Range("E30").Value = "  +4.56%  "
$ws.Range("D31").Value = "'34.14"
$ws.Range("E31").Value = "  +9.41%  "
$ws.Range("D32").Value = "'153.42"
$ws.Range("E32").Value = "  +2.78%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Value = "'4.93"
$ws.Range("E34").Value = "  +3.10%  "
$ws.Range("D35").Value = "'0.0721"
$ws.Range("E35").Value = "  +5.47%  "
$ws.Range("E36").Value = "  +2.35%  "
$ws.Range("D37").Value = "'16.78"
$ws.Range("E37").Value = "  +14.75%  "
$ws.Range("E38").Value = "  +6.84%  "
$ws.Range("E39").Value = "  +2.68%  "
$ws.Range("D40").Value = "'2.74"
$ws.Range("E40").Value = "  +2.38%  "
$ws.Range("E41").Value = "  +5.94%  "
$ws.Range("E42").Value = "  +5.76%  "
$ws.Range("D43").Value = "2.009.32"
$ws.Range("E43").Value = "  +5.75%  "
$ws.Range("E44").Value = "  +7.20%  "
$ws.Range("E45").Value = "  +7.35%  "
$ws.Range("D46").Value = "'10.09"
$ws.Range("E46").Value = "  +12.14%  "
$ws.Range("D47").Value = "'16.42"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("E48").Value = "  +2.88%  "
$ws.Range("D49").Value = "2.464.31"
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("D50").Value = "'71.76"
$ws.Range("E50").Value = "  +7.09%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'90.30"
$ws.Range("E51").Value = "  +3.12%  "
